$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "DIAL IN ACCESS"
$ws.Range("D1").Value = "ENABLE"
$ws.Range("C2").Value = "DIAL IN ACCESS"
$ws.Range("D2").Value = "DISABLE"
$ws.Range("C3").Value = "I87 LBSW"
$ws.Range("C4").Value = "I87 LBSW"
$ws.Range("C5").Value = "NO 6 LBSW"
$ws.Range("C6").Value = "NO 6 LBSW"
$ws.Range("C7").Value = "A53 LBSW"
$ws.Range("C8").Value = "A53 LBSW"
$ws.Range("D8").Value = " CLOSE"
$ws.Range("C9").Value = "BK 1 VCB"
$ws.Range("C10").Value = "BK 1 VCB"
$ws.Range("C11").Value = "47417 VCB"
$ws.Range("D11").Value = "OPEN"
$ws.Range("C12").Value = "47417 VCB"
$ws.Range("D12").Value = "CLOSE"
$ws.Range("C13").Value = "BT VCB"
$ws.Range("D13").Value = "OPEN"
$ws.Range("C14").Value = "BT VCB"
$ws.Range("D14").Value = "CLOSE"
$ws.Range("C15").Value = "BK 2 VCB"
$ws.Range("D15").Value = "OPEN"
$ws.Range("C16").Value = "BK 2 VCB"
$ws.Range("D16").Value = "CLOSE"
$ws.Range("C17").Value = "47418 VCB"
$ws.Range("D17").Value = "OPEN"
$ws.Range("C18").Value = "47418 VCB"
$ws.Range("D18").Value = "CLOSE"
$ws.Range("C19").Value = "47419 VCB"
$ws.Range("D19").Value = "OPEN"
$ws.Range("C20").Value = "47419 VCB"
$ws.Range("D20").Value = "CLOSE"
$ws.Range("C21").Value = "47416 VCB"
$ws.Range("D21").Value = "OPEN"
$ws.Range("C22").Value = "47416 VCB"
$ws.Range("D22").Value = "CLOSE"
$ws.Range("C23").Value = "CAP 1 VCB"
$ws.Range("D23").Value = "OPEN"
$ws.Range("C24").Value = "CAP 1 VCB"
$ws.Range("D24").Value = "CLOSE"
$ws.Range("C25").Value = "47415 VCB"
$ws.Range("D25").Value = "OPEN"
$ws.Range("C26").Value = "47415 VCB"
$ws.Range("D26").Value = "CLOSE"
$ws.Range("C27").Value = "BK 1 CKT INT CI-2"
$ws.Range("D27").Value = "OPEN"
$ws.Range("C28").Value = "BK 1 CKT INT CI-2"
$ws.Range("D28").Value = "CLOSE"
$ws.Range("C29").Value = "ADAPTIVE RELAYINIG"
$ws.Range("D29").Value = "OFF"
$ws.Range("C30").Value = "ADAPTIVE RELAYINIG"
$ws.Range("D30").Value = "ON"
$ws.Range("C31").Value = "UNDEFINED"
$ws.Range("C32").Value = "UNDEFINED"

$ws.Range("G9").Select()
